$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post about "あまり考えるな。疲れるよ" (row 309) was removed.
# Deleting the entire row shifts all subsequent rows up by one,
# matching the renumbering seen across the rest of the sheet.
$ws.Rows.Item(309).Delete()
